$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row ----
$ws.Range("A1").Value = "Task"
$ws.Range("B1").Value = "Start_Date"
$ws.Range("C1").Value = "End_Date"
$ws.Range("D1").Value = "Status"
$ws.Range("E1").Value = "Progress"

# ---- Row 2: Requirement Analysis ----
$ws.Range("A2").Value = "Requirement Analysis"
$ws.Range("B2").Value = 45200
$ws.Range("C2").Value = 45208
$ws.Range("D2").Value = "Completed"
$ws.Range("E2").Value = 1

# ---- Row 3: System Design ----
$ws.Range("A3").Value = "System Design"
$ws.Range("B3").Value = 45209
$ws.Range("C3").Value = 45223
$ws.Range("D3").Value = "Completed"
$ws.Range("E3").Value = 1

# ---- Row 4: Implementation ----
$ws.Range("A4").Value = "Implementation"
$ws.Range("B4").Value = 45224
$ws.Range("C4").Value = 45244
$ws.Range("D4").Value = "In Progress"
$ws.Range("E4").Value = 0.6

# ---- Row 5: Testing ----
$ws.Range("A5").Value = "Testing"
$ws.Range("B5").Value = 45245
$ws.Range("C5").Value = 45254
$ws.Range("D5").Value = "Pending"
$ws.Range("E5").Value = 0

# ---- Row 6: Deployment ----
$ws.Range("A6").Value = "Deployment"
$ws.Range("B6").Value = 45255
$ws.Range("C6").Value = 45260
$ws.Range("D6").Value = "Pending"
$ws.Range("E6").Value = 0

# ---- Date formatting for Start_Date / End_Date columns ----
$ws.Range("B2:C6").NumberFormat = "mm-dd-yy"

# ---- Column widths (approximate "Task/Dates/Status/Progress" layout) ----
$ws.Columns.Item(1).ColumnWidth = 25.71
$ws.Columns.Item(2).ColumnWidth = 15.71
$ws.Columns.Item(3).ColumnWidth = 15.71
$ws.Columns.Item(4).ColumnWidth = 15.71
$ws.Columns.Item(5).ColumnWidth = 10.71

# ---- Page setup ----
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
